$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 13 & 14 previously had no entries in columns B ("Tareas Completadas
# en la Semana") and C ("Tareas Restantes (Real)"). Copy the formatting
# from row 12 (which already carries the shared style) onto B13:C14 before
# filling in the new values so the new cells pick up the same cell style.
$ws.Range("B12:C12").Copy()
$ws.Range("B13:C14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Updated burndown data (columns B, C, D) for rows 2-14.
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 98
$ws.Range("D2").Value = 98

$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 95
$ws.Range("D3").Value = 90

$ws.Range("B4").Value = 8
$ws.Range("C4").Value = 87
$ws.Range("D4").Value = 82

$ws.Range("B5").Value = 15
$ws.Range("C5").Value = 72
$ws.Range("D5").Value = 74

$ws.Range("B6").Value = 6
$ws.Range("C6").Value = 66
$ws.Range("D6").Value = 66

$ws.Range("B7").Value = 6
$ws.Range("C7").Value = 60
$ws.Range("D7").Value = 58

$ws.Range("B8").Value = 10
$ws.Range("C8").Value = 50
$ws.Range("D8").Value = 50

$ws.Range("B9").Value = 12
$ws.Range("C9").Value = 38
$ws.Range("D9").Value = 42

$ws.Range("B10").Value = 5
$ws.Range("C10").Value = 33
$ws.Range("D10").Value = 34

$ws.Range("B11").Value = 7
$ws.Range("C11").Value = 26
$ws.Range("D11").Value = 26

$ws.Range("B12").Value = 7
$ws.Range("C12").Value = 19
$ws.Range("D12").Value = 18

$ws.Range("B13").Value = 7
$ws.Range("C13").Value = 12
$ws.Range("D13").Value = 10

$ws.Range("B14").Value = 12
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0
